$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - Ar da sala / A/C / 23 / FALSE
$ws.Cells.Item(1, 1).Value = "Ar da sala"
$ws.Cells.Item(1, 2).Value = "A/C"
$ws.Cells.Item(1, 3).Value = 23
$ws.Cells.Item(1, 4).Value = $false

# Row 2 - Ar da cozinha / A/C / 23 / FALSE
$ws.Cells.Item(2, 1).Value = "Ar da cozinha"
$ws.Cells.Item(2, 2).Value = "A/C"
$ws.Cells.Item(2, 3).Value = 23
$ws.Cells.Item(2, 4).Value = $false

# Row 3 - Ar do quarto / A/C / 23 / FALSE
$ws.Cells.Item(3, 1).Value = "Ar do quarto"
$ws.Cells.Item(3, 2).Value = "A/C"
$ws.Cells.Item(3, 3).Value = 23
$ws.Cells.Item(3, 4).Value = $false

# Row 4 - a / Lâmpada / 0 / FALSE
$ws.Cells.Item(4, 1).Value = "a"
$ws.Cells.Item(4, 2).Value = "Lâmpada"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = $false

# Row 5 - tv da sala / Televisor / 1 / 0 / FALSE
$ws.Cells.Item(5, 1).Value = "tv da sala"
$ws.Cells.Item(5, 2).Value = "Televisor"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = $false
